# Adds new survey submission rows (2024-12-06) that arrived after the
# workbook was last exported, across the relevant sheets of the
# "plastiq_input_information" workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# contact_data
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("contact_data")

$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "2024-12-06 11:06:02"
$ws.Range("E12").Value = "SKZ"

$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "2024-12-06 13:41:46"

$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "2024-12-06 13:42:27"
$ws.Range("C14").Value = "Stefan"
$ws.Range("D14").Value = "Trieß"
$ws.Range("E14").Value = "SKZ"
$ws.Range("F14").Value = "Engineer"
$ws.Range("G14").Value = "s.triess@skz.de"
$ws.Range("H14").Value = "+49 15786027870"

$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "2024-12-06 15:39:33"
$ws.Range("E15").Value = "SKZ"

# ---------------------------------------------------------------------
# company_data
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("company_data")

$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "2024-12-06 11:06:23"
$ws.Range("C13").Value = "SKZ"
$ws.Range("D13").Value = "Friedrich-Bergius-Ring 22"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "97076"
$ws.Range("F13").Value = "Würzburg"
$ws.Range("G13").Value = "Bayern"
$ws.Range("H13").Value = "Deutschland"
$ws.Range("I13").Value = $False
$ws.Range("J13").Value = $False
$ws.Range("L13").Value = 49.80282025
$ws.Range("M13").Value = 10.00010726291456

$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "2024-12-06 13:43:24"
$ws.Range("C14").Value = "SKZ"
$ws.Range("D14").Value = "Friedrich-Bergius-Ring 22"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "97076"
$ws.Range("F14").Value = "Würzburg"
$ws.Range("G14").Value = "Bayern"
$ws.Range("H14").Value = "Deutschland"
$ws.Range("I14").Value = $True
$ws.Range("J14").Value = $False
$ws.Range("L14").Value = 49.80282025
$ws.Range("M14").Value = 10.00010726291456

$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "2024-12-06 14:42:20"
$ws.Range("D15").Value = "Friedrich Bergius Ring 22"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "97076"
$ws.Range("F15").Value = "Würzburg"
$ws.Range("G15").Value = "Bayern"
$ws.Range("H15").Value = "Deutschland"
$ws.Range("I15").Value = $False
$ws.Range("J15").Value = $False
$ws.Range("L15").Value = 49.80282025
$ws.Range("M15").Value = 10.00010726291456

$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "2024-12-06 14:42:22"
$ws.Range("D16").Value = "Friedrich Bergius Ring 22"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "97076"
$ws.Range("F16").Value = "Würzburg"
$ws.Range("G16").Value = "Bayern"
$ws.Range("H16").Value = "Deutschland"
$ws.Range("I16").Value = $False
$ws.Range("J16").Value = $False
$ws.Range("L16").Value = 49.80282025
$ws.Range("M16").Value = 10.00010726291456

$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "2024-12-06 15:40:03"
$ws.Range("C17").Value = "SKZ"
$ws.Range("D17").Value = "Friedrich-Bergius-Ring 22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "97076"
$ws.Range("F17").Value = "Würzburg"
$ws.Range("G17").Value = "Bayern"
$ws.Range("H17").Value = "Deutschland"
$ws.Range("I17").Value = $True
$ws.Range("J17").Value = $False
$ws.Range("L17").Value = 49.80282025
$ws.Range("M17").Value = 10.00010726291456

# ---------------------------------------------------------------------
# product_fractions
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("product_fractions")

$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "2024-12-06 11:06:39"
$ws.Range("C12").Value = "['PE-LLD', 'PS', 'Magnesium']"
$ws.Range("D12").Value = "['', '', '']"
$ws.Range("E12").Value = "[80.0, 10.0, 10.0]"

$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "2024-12-06 13:46:50"
$ws.Range("C13").Value = "['PP', 'PS', 'PVC-U']"
$ws.Range("D13").Value = "['', '', '']"
$ws.Range("E13").Value = "[40.0, 30.0, 30.0]"

$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "2024-12-06 14:43:02"
$ws.Range("C14").Value = "['PE-LD', 'PS', 'Duromere', 'PUR']"
$ws.Range("D14").Value = "['', '', '', '']"
$ws.Range("E14").Value = "[40.0, 30.0, 10.0, 20.0]"

$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "2024-12-06 15:40:29"
$ws.Range("C15").Value = "['PP', 'ABS', 'PMMA', 'PUR']"
$ws.Range("D15").Value = "['', '', '', '']"
$ws.Range("E15").Value = "[40.0, 20.0, 20.0, 20.0]"

# ---------------------------------------------------------------------
# product_origin
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("product_origin")

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "2024-12-06 11:07:16"
$ws.Range("C8").Value = "Post-Consumer (PC) – getrennte Sammlung"
$ws.Range("D8").Value = "Kittel Reinraum"
$ws.Range("E8").Value = "lokal als Bringsystem"

# ---------------------------------------------------------------------
# product_quality
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("product_quality")

$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "2024-12-06 11:07:55"
$ws.Range("C10").Value = "Ja"
$ws.Range("D10").Value = "divers"
$ws.Range("E10").Value = 100
$ws.Range("F10").Value = "keine"
$ws.Range("H10").Value = "[]"
$ws.Range("I10").Value = "[]"

$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "2024-12-06 13:47:16"
$ws.Range("C11").Value = "Ja"
$ws.Range("D11").Value = "weiß"
$ws.Range("E11").Value = 99.98999999999999
$ws.Range("F11").Value = "gering"
$ws.Range("H11").Value = "[[], [], []]"
$ws.Range("I11").Value = "[[], [], []]"

$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "2024-12-06 14:43:43"
$ws.Range("C12").Value = "Ja"
$ws.Range("D12").Value = "natur"
$ws.Range("E12").Value = 100
$ws.Range("F12").Value = "hoch"
$ws.Range("H12").Value = "[]"
$ws.Range("I12").Value = "[]"

# ---------------------------------------------------------------------
# product_amount
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("product_amount")

$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "2024-12-06 11:09:39"
$ws.Range("C13").Value = 10
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = "Woche"

$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "2024-12-06 13:47:34"
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = "Quartal"

$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "2024-12-06 14:44:07"
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = "Quartal"

$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "2024-12-06 15:41:25"
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = "Monat"
